{"js": "// Resume updates:\n//  1. \"6-7 months\" -> \"7-9 months\" (MDN Web Docs bullet)\n//  2. Remove \"C#, \" from the languages/technologies list\n//  3. \"Visual Studio, Git, Eclipse\" -> \"Windows, Visual Studio, Unix, Linux, Git, Eclipse\"\n\nconst body = context.document.body;\n\n// 1) Update the \"6-7 months\" -> \"7-9 months\" duration.\nconst monthsHits = body.search(\"over the course of 6-7 months\", { matchCase: true });\nmonthsHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < monthsHits.items.length; i++) {\n  monthsHits.items[i].insertText(\"over the course of 7-9 months\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Drop \"C#\" from the C++/Java/... skills line.\nconst skillsHits = body.search(\"Python, C#, SQL\", { matchCase: true });\nskillsHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < skillsHits.items.length; i++) {\n  skillsHits.items[i].insertText(\"Python, SQL\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Expand the tools line with Windows / Unix / Linux.\nconst toolsHits = body.search(\"Visual Studio, Git, Eclipse\", { matchCase: true });\ntoolsHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < toolsHits.items.length; i++) {\n  toolsHits.items[i].insertText(\"Windows, Visual Studio, Unix, Linux, Git, Eclipse\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Resume updates:\n#  1. \"6-7 months\" -> \"7-9 months\" (MDN Web Docs bullet)\n#  2. Remove \"C#, \" from the languages/technologies list\n#  3. \"Visual Studio, Git, Eclipse\" -> \"Windows, Visual Studio, Unix, Linux, Git, Eclipse\"\n\n$d = $word.ActiveDocument\n\n# 1) Update the \"6-7 months\" -> \"7-9 months\" duration.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"over the course of 6-7 months\"\n$find1.Replacement.Text = \"over the course of 7-9 months\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) Drop \"C#\" from the C++/Java/... skills line.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Python, C#, SQL\"\n$find2.Replacement.Text = \"Python, SQL\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# 3) Expand the tools line with Windows / Unix / Linux.\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"Visual Studio, Git, Eclipse\"\n$find3.Replacement.Text = \"Windows, Visual Studio, Unix, Linux, Git, Eclipse\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n"}
